$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value2 = 6.41
$ws.Range("C5").Value2 = 5.12
$ws.Range("D5").Value2 = 0.08
$ws.Range("E5").Value2 = 12.88
$ws.Range("F5").Value2 = 11.82
$ws.Range("G5").Value2 = 4.63
$ws.Range("H5").Value2 = 22.22
$ws.Range("I5").Value2 = 7.27
$ws.Range("J5").Value2 = 3.38
$ws.Range("K5").Value2 = 5.39
$ws.Range("L5").Value2 = 5.68
$ws.Range("M5").Value2 = 5.96
$ws.Range("N5").Value2 = 1.54
$ws.Range("O5").Value2 = 4.89
$ws.Range("P5").Value2 = 6.73
$ws.Range("Q5").Value2 = 3.68
$ws.Range("R5").Value2 = 0.37
$ws.Range("S5").Value2 = 0.32
$ws.Range("T5").Value2 = 66.90000000000001
$ws.Range("U5").Value2 = 13.49
$ws.Range("V5").Value2 = 4.72
$ws.Range("W5").Value2 = 9.58
$ws.Range("X5").Value2 = 4.84
$ws.Range("Y5").Value2 = 0.64
$ws.Range("Z5").Value2 = 9.93
$ws.Range("AA5").Value2 = 3.67
$ws.Range("AB5").Value2 = 3.62
$ws.Range("AC5").Value2 = 4.23
$ws.Range("AD5").Value2 = 6.09
$ws.Range("AE5").Value2 = 0
$ws.Range("AF5").Value2 = 19.72
$ws.Range("AG5").Value2 = 2.71
$ws.Range("AH5").Value2 = 5.48

$ws.Rows.Item(6).Delete()

